# Generate Report for Archive
#
# 1. The shared string "Ready for handoff" becomes "In Translation" -
#    update every cell that shows that status (Overview!E2:F3, and the
#    "Status" column [C] on the zh-cn and de-de detail sheets).
# 2. Because the new status text is shorter, the "Status" columns are
#    narrowed (from ~17.22 chars down to ~13.41 chars of raw width).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text wherever "Ready for handoff" appears ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- 2. Narrow the now-shorter "Status" columns ---
# ColumnWidth is in characters and the host snaps it to the same
# pixel grid Excel itself uses, so 12.5 is the input that lands on the
# closest representable width to the target ~13.41.
$overview.Range("E1").EntireColumn.ColumnWidth = 12.5
$overview.Range("F1").EntireColumn.ColumnWidth = 12.5

$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
